# The workbook contains one sheet with a weekly price list for Mango at
# "Terminal La Palmera de La Serena". A new week of data (rows for
# Calidad = Especial / Primera / Segunda, with Fecha = 44452) needs to be
# inserted right above the current row 249, pushing every following row
# (249-380) down by three rows (they become 252-383, unchanged).
#
# The new week's values mirror the most recent existing week in the sheet
# (the old rows 373-375, i.e. the rows that end up at 376-378 after the
# shift) but dated 4 days later (44448 + 4 = 44452).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert three blank rows above row 249; this shifts the old rows
#    249:380 down to 252:383 and grows the sheet to A1:T383.
$ws.Rows("249:251").Insert()

# 2) The block that used to be rows 373:375 now lives at 376:378 (it was
#    also shifted down by the insert). Copy it into the freshly inserted
#    rows 249:251 as the template for the new week.
$ws.Range("A376:T378").Copy()
$ws.Range("A249").PasteSpecial()

# 3) Update the date for the new week (one week/4 days after 44448).
$ws.Range("D249:D251").Value = 44452
